$d = $word.ActiveDocument

# 1. Replace "nonmatriculated" with "non-matriculated"
$d.Content.Find.Execute("nonmatriculated", $true, $false, $false, $false, $false,
                         $true, 1, $false, "non-matriculated", 2)

# 2. Set the page orientation explicitly to portrait (adds w:orient="portrait" to pgSz)
$d.PageSetup.Orientation = 0
